$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '25.698.44'
Set-TextValue 'E2' '  -3.75%  '
Set-TextValue 'D3' '1.745.50'
Set-TextValue 'E3' '  -5.71%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '236.29'
Set-TextValue 'E5' '  -9.95%  '
Set-TextValue 'E6' '  +0.07%  '
Set-TextValue 'D7' '0.4922'
Set-TextValue 'E7' '  -8.05%  '
Set-TextValue 'D8' '41.62'
Set-TextValue 'E8' '  -8.01%  '
Set-TextValue 'D9' '0.2508'
Set-TextValue 'E9' '  -21.36%  '
Set-TextValue 'D10' '0.06005'
Set-TextValue 'E10' '  -14.00%  '
Set-TextValue 'D11' '1.744.06'
Set-TextValue 'E11' '  -5.78%  '
Set-TextValue 'D12' '0.06838'
Set-TextValue 'E12' '  -12.70%  '
Set-TextValue 'D13' '14.81'
Set-TextValue 'E13' '  -22.04%  '
Set-TextValue 'D14' '4.443'
Set-TextValue 'E14' '  -12.20%  '
Set-TextValue 'D15' '76.83'
Set-TextValue 'E15' '  -14.46%  '
Set-TextValue 'E16' '  -27.31%  '
Set-TextValue 'D17' '1.001'
Set-TextValue 'E17' '  -0.01%  '
Set-TextValue 'D18' '1.001'
Set-TextValue 'E18' '  +0.05%  '
Set-TextValue 'D19' '25.738.54'
Set-TextValue 'E19' '  -3.71%  '
Set-TextValue 'E20' '  -20.66%  '
Set-TextValue 'D21' '0.000006547'
Set-TextValue 'E21' '  -18.30%  '
Set-TextValue 'D22' '1.965.32'
Set-TextValue 'E22' '  -5.59%  '
Set-TextValue 'D23' '3.990'
Set-TextValue 'E23' '  -14.47%  '
Set-TextValue 'D24' '5.002'
Set-TextValue 'E24' '  -17.33%  '
Set-TextValue 'D25' '7.856'
Set-TextValue 'E25' '  -16.55%  '
Set-TextValue 'D26' '136.68'
Set-TextValue 'E26' '  -4.42%  '
Set-TextValue 'E27' '  -12.79%  '
Set-TextValue 'D28' '1.802'
Set-TextValue 'E28' '  -18.78%  '
Set-TextValue 'D29' '14.63'
Set-TextValue 'E29' '  -14.71%  '
Set-TextValue 'D30' '101.81'
Set-TextValue 'E30' '  -8.98%  '
Set-TextValue 'D31' '3.755'
Set-TextValue 'E31' '  -13.44%  '
Set-TextValue 'D32' '0.07982'
Set-TextValue 'E32' '  -8.91%  '
Set-TextValue 'D33' '3.366'
Set-TextValue 'E33' '  -18.31%  '
Set-TextValue 'D34' '0.04385'
Set-TextValue 'E34' '  -10.11%  '
Set-TextValue 'D35' '0.9999'
Set-TextValue 'E35' '  -0.01%  '
Set-TextValue 'D36' '2.637'
Set-TextValue 'E36' '  -9.03%  '
Set-TextValue 'D37' '0.9677'
Set-TextValue 'E37' '  -15.40%  '
Set-TextValue 'D38' '0.6015'
Set-TextValue 'E38' '  -19.01%  '
Set-TextValue 'D39' '2.679'
Set-TextValue 'E39' '  -14.01%  '
Set-TextValue 'D40' '2.000'
Set-TextValue 'E40' '  -15.49%  '
Set-TextValue 'E41' '  +0.08%  '
Set-TextValue 'D42' '102.63'
Set-TextValue 'E42' '  -6.13%  '
Set-TextValue 'D43' '0.01500'
Set-TextValue 'E43' '  -14.44%  '
Set-TextValue 'D44' '0.7544'
Set-TextValue 'E44' '  -16.86%  '
Set-TextValue 'D45' '5.156'
Set-TextValue 'D46' '0.3715'
Set-TextValue 'E46' '  -23.39%  '
Set-TextValue 'D47' '0.05258'
Set-TextValue 'E47' '  -9.92%  '
Set-TextValue 'D48' '0.1063'
Set-TextValue 'E48' '  -15.32%  '
Set-TextValue 'D49' '29.99'
Set-TextValue 'E49' '  -14.87%  '
Set-TextValue 'D50' '5.865'
Set-TextValue 'E50' '  -24.17%  '
Set-TextValue 'E51' '  -13.91%  '
